$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: update the rich-text "Volume 31   Number  13" -> "...14"
# and the reporting week date range, editing only the sub-run(s) of text
# that changed (processed back-to-front so earlier offsets stay valid) and
# re-applying the original run font so the rich text formatting is kept.
# ---------------------------------------------------------------------------
$a8 = $ws.Range("A8")
$a8r4 = $a8.Characters(21, 2)
$a8r4.Text = "14"
$a8r4.Font.Size = 10
$a8r4.Font.Name = "Andale WT"
$a8r3 = $a8.Characters(10, 11)
$a8r3.Font.Size = 10
$a8r3.Font.Name = "Andale WT"
$a8r2 = $a8.Characters(8, 2)
$a8r2.Font.Size = 10
$a8r2.Font.Name = "Andale WT"
$a8r1 = $a8.Characters(1, 7)
$a8r1.Font.Size = 10
$a8r1.Font.Name = "Andale WT"

$c9 = $ws.Range("C9")
$c9r4 = $c9.Characters(47, 9)
$c9r4.Text = "4/7/2024"
$c9r4.Font.Size = 10
$c9r4.Font.Name = "Andale WT"
$c9r3 = $c9.Characters(36, 11)
$c9r3.Font.Size = 10
$c9r3.Font.Name = "Andale WT"
$c9r2 = $c9.Characters(27, 9)
$c9r2.Text = "4/1/2024"
$c9r2.Font.Size = 10
$c9r2.Font.Name = "Andale WT"
$c9r1 = $c9.Characters(1, 26)
$c9r1.Font.Size = 10
$c9r1.Font.Name = "Andale WT"

# ---------------------------------------------------------------------------
# Weekly crime-stat numbers (rows 16-31). Plain numeric overwrites; Excel
# keeps the existing cell style/number format for each of these cells.
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 11
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 266.666666666667
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 64.285714285714
$ws.Range("I16").Value = 67
$ws.Range("J16").Value = 28
$ws.Range("K16").Value = 139.285714285714
$ws.Range("L16").Value = 81.081081081081
$ws.Range("M16").Value = 11.666666666666
$ws.Range("N16").Value = -78.246753246753
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 13.333333333333
$ws.Range("I17").Value = 72
$ws.Range("J17").Value = 68
$ws.Range("K17").Value = 5.882352941176
$ws.Range("L17").Value = 84.615384615384
$ws.Range("M17").Value = 125
$ws.Range("N17").Value = -32.075471698113
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -23.529411764705
$ws.Range("I18").Value = 48
$ws.Range("J18").Value = 63
$ws.Range("K18").Value = -23.809523809523
$ws.Range("L18").Value = -20
$ws.Range("M18").Value = 128.571428571429
$ws.Range("N18").Value = -74.054054054054
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = -9.090909090909
$ws.Range("I19").Value = 185
$ws.Range("J19").Value = 172
$ws.Range("K19").Value = 7.558139534883
$ws.Range("L19").Value = 14.906832298136
$ws.Range("M19").Value = 60.869565217391
$ws.Range("N19").Value = -30.188679245283
$ws.Range("C20").Value = 1
$ws.Range("I20").Value = 15
$ws.Range("K20").Value = 7.142857142857
$ws.Range("L20").Value = 15.384615384615
$ws.Range("M20").Value = -6.25
$ws.Range("N20").Value = -92.105263157894
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 58.333333333333
$ws.Range("G21").Value = 105
$ws.Range("H21").Value = 2.857142857142
$ws.Range("I21").Value = 388
$ws.Range("J21").Value = 348
$ws.Range("K21").Value = 11.494252873563
$ws.Range("L21").Value = 23.566878980891
$ws.Range("M21").Value = 58.367346938775
$ws.Range("N21").Value = -63.396226415094
$ws.Range("D22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 12
$ws.Range("K22").Value = 16.666666666666
$ws.Range("M22").Value = -26.315789473684
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 0
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = 9.677419354838
$ws.Range("F24").Value = 183
$ws.Range("G24").Value = 153
$ws.Range("H24").Value = 19.607843137254
$ws.Range("I24").Value = 575
$ws.Range("J24").Value = 505
$ws.Range("K24").Value = 13.861386138613
$ws.Range("L24").Value = 42.679900744416
$ws.Range("M24").Value = 64.756446991404
$ws.Range("C25").Value = 33
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = 26.923076923076
$ws.Range("F25").Value = 165
$ws.Range("G25").Value = 135
$ws.Range("H25").Value = 22.222222222222
$ws.Range("I25").Value = 521
$ws.Range("J25").Value = 459
$ws.Range("K25").Value = 13.507625272331
$ws.Range("L25").Value = 56.456456456456
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 39
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = 69.565217391304
$ws.Range("I26").Value = 126
$ws.Range("J26").Value = 99
$ws.Range("K26").Value = 27.272727272727
$ws.Range("L26").Value = 44.827586206896
$ws.Range("M26").Value = 17.757009345794
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("I28").Value = 9
$ws.Range("J28").Value = 13
$ws.Range("K28").Value = -30.769230769230
$ws.Range("N29").Value = -91.666666666666
$ws.Range("N30").Value = -90.909090909090
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 4
$ws.Range("K31").Value = -50
$ws.Range("L31").Value = 100

# ---------------------------------------------------------------------------
# A handful of cells flip from a numeric value to the special text markers
# used elsewhere in the sheet for "no data" ("0") and "not meaningful"
# ("***.*"). Setting .Value directly would just store a number (for "0")
# or create a brand-new shared string with a different style, so instead we
# write the text first (apostrophe-prefixed so "0" is kept as text) and
# then copy the number format/style from an existing matching cell.
# ---------------------------------------------------------------------------
$ws.Range("D20").Value = "'0"
$ws.Range("E20").Value = "***.*"
$ws.Range("C23").Value = "'0"
$ws.Range("D23").Value = "'0"
$ws.Range("E23").Value = "***.*"
$ws.Range("D31").Value = "'0"
$ws.Range("E31").Value = "***.*"

$ws.Range("C22").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").PasteSpecial(-4122)
$excel.CutCopyMode = 0
